$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NextBus1")

# EstimatedTimeOfArrival (F) and MinutesToArrival (U) were refreshed with a
# newer poll of the NextBus API; TypeOfBus (L) for row 11 also changed.
$ws.Range("F2").Value = 45684.98508101852
$ws.Range("U2").Value = 1

$ws.Range("F3").Value = 45684.98773148148
$ws.Range("U3").Value = 5

$ws.Range("F4").Value = 45684.99252314815
$ws.Range("U4").Value = 12

$ws.Range("F5").Value = 45684.98704861111
$ws.Range("U5").Value = 4

$ws.Range("F6").Value = 45684.99150462963
$ws.Range("U6").Value = 10

$ws.Range("F7").Value = 45684.98303240741
$ws.Range("U7").Value = -1

$ws.Range("F8").Value = 45684.98582175926
$ws.Range("U8").Value = 2

$ws.Range("F9").Value = 45684.98793981481
$ws.Range("U9").Value = 5

$ws.Range("F10").Value = 45684.9884375
$ws.Range("U10").Value = 6

$ws.Range("F11").Value = 45684.99438657407
$ws.Range("L11").Value = "SD"
$ws.Range("U11").Value = 14

$ws.Range("F12").Value = 45684.9946875
$ws.Range("U12").Value = 15

$ws.Range("F13").Value = 45684.98508101852
$ws.Range("U13").Value = 1
